# This document renders its visible content almost entirely through a very
# large number of absolutely-positioned text boxes (mc:AlternateContent /
# wp:anchor drawings) that all live inside a single run of a single
# paragraph. Because of that, none of the date text is reachable through
# Document.Content.Text / Find.Execute (those only see the "flow" text,
# which is empty here), and the Shapes collection's write path in this
# runtime does not reliably address individual shapes by index.
#
# The reliable way to edit this content is to pull the whole document's
# OOXML via Range.WordOpenXML, perform precise, uniquely-scoped textual
# substitutions that mirror the target unified diff, and push the result
# back with Range.InsertXML (which replaces the contents of the range).

$d = $word.ActiveDocument
$r = $d.Content
$xml = $r.WordOpenXML

# --- Shape99 ("2023" label next to "JUNE 29,") -----------------------------
# wp:anchor positionH posOffset moves right because the neighbouring
# "JUNE 29," box (Shape98) grew wider.
$xml = $xml.Replace(
  "<wp:posOffset>880745</wp:posOffset>",
  "<wp:posOffset>1096645</wp:posOffset>")

# mc:Fallback VML shape margin-left shifts accordingly (width unchanged).
$xml = $xml.Replace(
  "margin-left:69.35pt;margin-top:807.1pt;width:20.3pt;height:9.3pt",
  "margin-left:86.35pt;margin-top:807.1pt;width:20.3pt;height:9.3pt")

# --- Shape98 ("JUNE 29," uppercase label) ----------------------------------
# wp:extent / wp:effectExtent grow to fit the new, longer text.
$xml = $xml.Replace(
  "<wp:extent cx=`"436880`" cy=`"118745`"/><wp:effectExtent l=`"5715`" t=`"0`" r=`"0`" b=`"0`"/>",
  "<wp:extent cx=`"734060`" cy=`"118745`"/><wp:effectExtent l=`"0`" t=`"0`" r=`"0`" b=`"0`"/>")

# DrawingML shape's own extent (a:xfrm/a:ext) matches the new width.
$xml = $xml.Replace(
  "<a:ext cx=`"437040`" cy=`"118800`"/>",
  "<a:ext cx=`"734040`" cy=`"118800`"/>")

# mc:Fallback VML shape width grows too (position unchanged).
$xml = $xml.Replace(
  "margin-left:39.65pt;margin-top:807.1pt;width:34.35pt;height:9.3pt",
  "margin-left:39.65pt;margin-top:807.1pt;width:57.75pt;height:9.3pt")

# Text itself: both the mc:Choice (DrawingML) and mc:Fallback (VML) runs
# carry their own copy of the text - replace every occurrence.
$xml = $xml.Replace("<w:t>JUNE 29,</w:t>", "<w:t>SEPTEMBER 7,</w:t>")

# --- Shape24 ("June 29," title-case label, upper-right area) --------------
# wp:anchor positionH posOffset moves left because the box grew wider and
# stays right-aligned against its neighbour.
$xml = $xml.Replace(
  "<wp:posOffset>6410960</wp:posOffset>",
  "<wp:posOffset>6154420</wp:posOffset>")

# wp:extent / wp:effectExtent grow to fit the new, longer text.
$xml = $xml.Replace(
  "<wp:extent cx=`"465455`" cy=`"132715`"/><wp:effectExtent l=`"19050`" t=`"0`" r=`"0`" b=`"0`"/>",
  "<wp:extent cx=`"774065`" cy=`"132715`"/><wp:effectExtent l=`"0`" t=`"0`" r=`"0`" b=`"0`"/>")

# DrawingML shape's own extent (a:xfrm/a:ext) matches the new width.
$xml = $xml.Replace(
  "<a:ext cx=`"465480`" cy=`"132840`"/>",
  "<a:ext cx=`"774000`" cy=`"132840`"/>")

# mc:Fallback VML shape shifts left and grows wider.
$xml = $xml.Replace(
  "margin-left:504.8pt;margin-top:163.7pt;width:36.6pt;height:10.4pt",
  "margin-left:484.6pt;margin-top:163.7pt;width:60.9pt;height:10.4pt")

# Text itself: both the mc:Choice (DrawingML) and mc:Fallback (VML) runs
# carry their own copy of the text - replace every occurrence.
$xml = $xml.Replace("<w:t>June 29,</w:t>", "<w:t>September 7,</w:t>")

$r.InsertXML($xml)
